$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for the new columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style from an existing header cell (e.g. AC1) to the new headers
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null # xlPasteFormats

# Fill season record data for each team row (rows 2-40)
for ($r = 2; $r -le 40; $r++) {
    $ws.Cells.Item($r, 30).Value = 75   # AD = column 30, Wins
    $ws.Cells.Item($r, 31).Value = 87   # AE = column 31, Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF = column 32, Ties
}
